# Fixed some f12018 motion data mappings
$wb = $excel.ActiveWorkbook

# --- CarData sheet: fill in the G column (running totals) for rows 5-16 ---
$wsCar = $wb.Worksheets.Item("CarData")
$wsCar.Range("G5:G16").Formula = "=G4+F4"

# --- Selections / active sheet tab ---
$wsMotion = $wb.Worksheets.Item("Motion")
$wsMotion.Range("C5").Select()

$wsCar.Activate()
$wsCar.Range("A1:C1").Select()
